$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Species" cell (F) already contains an all-caps species name;
# these just need to be re-cased to title case. Did Harvest Occur? (B) and
# Unknown Sex Count (J) are already correct on these rows.
$speciesFix = @{ 12 = "Weasel"; 27 = "Lynx"; 33 = "Marten"; 38 = "Mink" }

for ($row = 2; $row -le 88; $row++) {
    if ($speciesFix.ContainsKey($row)) {
        $ws.Cells.Item($row, 6).Value = $speciesFix[$row]
    } else {
        $ws.Cells.Item($row, 2).Value = "Yes"
        $ws.Cells.Item($row, 6).Value = "Na"
        $ws.Cells.Item($row, 10).Value = 1
    }
}
